$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "22.020.03"
Set-TextValue "E2" "  +6.96%  "
Set-TextValue "D3" "1.579.89"
Set-TextValue "E3" "  +6.75%  "
Set-TextValue "D4" "0.9993"
Set-TextValue "E4" "  -0.58%  "
Set-TextValue "D5" "0.9854"
Set-TextValue "E5" "  +1.79%  "
Set-TextValue "D6" "287.86"
Set-TextValue "E6" "  +3.83%  "
Set-TextValue "D7" "0.3708"
Set-TextValue "E7" "  +1.33%  "
Set-TextValue "D8" "0.3296"
Set-TextValue "E8" "  +7.52%  "
Set-TextValue "D9" "42.77"
Set-TextValue "E9" "  +5.22%  "
Set-TextValue "D10" "1.152"
Set-TextValue "E10" "  +8.36%  "
Set-TextValue "D11" "0.07061"
Set-TextValue "E11" "  +6.18%  "
Set-TextValue "D12" "0.9949"
Set-TextValue "E12" "  -0.20%  "
Set-TextValue "D13" "20.22"
Set-TextValue "E13" "  +11.07%  "
Set-TextValue "D14" "5.881"
Set-TextValue "E14" "  +7.00%  "
Set-TextValue "D15" "6.552"
Set-TextValue "E15" "  +5.77%  "
Set-TextValue "D16" "0.9856"
Set-TextValue "E16" "  +1.68%  "
Set-TextValue "D17" "0.00001077"
Set-TextValue "E17" "  +4.00%  "
Set-TextValue "D18" "1.576.73"
Set-TextValue "E18" "  +6.60%  "
Set-TextValue "D19" "0.06437"
Set-TextValue "E19" "  +8.69%  "
Set-TextValue "D20" "75.67"
Set-TextValue "E20" "  +8.70%  "
Set-TextValue "E21" "  +10.77%  "
Set-TextValue "D22" "5.890"
Set-TextValue "E22" "  +7.83%  "
Set-TextValue "D23" "11.71"
Set-TextValue "E23" "  +5.61%  "
Set-TextValue "D24" "21.965.54"
Set-TextValue "E24" "  +6.56%  "
Set-TextValue "D25" "2.358"
Set-TextValue "E25" "  +4.44%  "
Set-TextValue "D26" "2.430"
Set-TextValue "E26" "  +13.20%  "
Set-TextValue "D27" "149.98"
Set-TextValue "E27" "  +6.37%  "
Set-TextValue "D28" "18.75"
Set-TextValue "E28" "  +8.23%  "
Set-TextValue "D29" "1.747.53"
Set-TextValue "E29" "  +6.86%  "
Set-TextValue "D30" "120.36"
Set-TextValue "E30" "  +5.49%  "
Set-TextValue "D31" "4.178"
Set-TextValue "E31" "  +6.09%  "
Set-TextValue "B32" "Filecoin"
Set-TextValue "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "5.689"
Set-TextValue "E32" "  +14.16%  "
Set-TextValue "B33" "ImmutableX"
Set-TextValue "C33" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D33" "0.9301"
Set-TextValue "E33" "  +13.78%  "
Set-TextValue "D34" "0.08243"
Set-TextValue "E34" "  +4.00%  "
Set-TextValue "D35" "1.655"
Set-TextValue "E35" "  +6.44%  "
Set-TextValue "B36" "InternetComputer(DFINITY)"
Set-TextValue "C36" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D36" "5.216"
Set-TextValue "E36" "  +10.04%  "
Set-TextValue "D37" "11.86"
Set-TextValue "E37" "  +13.04%  "
Set-TextValue "B38" "FraxShare"
Set-TextValue "C38" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D38" "8.626"
Set-TextValue "E38" "  +12.61%  "
Set-TextValue "D39" "0.06176"
Set-TextValue "E39" "  +5.66%  "
Set-TextValue "D40" "1.242"
Set-TextValue "E40" "  +2.26%  "
Set-TextValue "D41" "0.02189"
Set-TextValue "E41" "  +7.11%  "
Set-TextValue "D42" "0.2013"
Set-TextValue "E42" "  +6.85%  "
Set-TextValue "D43" "0.9857"
Set-TextValue "E43" "  +1.74%  "
Set-TextValue "D44" "0.5842"
Set-TextValue "E44" "  +10.01%  "
Set-TextValue "D45" "13.00"
Set-TextValue "E45" "  +6.57%  "
Set-TextValue "D46" "3.653"
Set-TextValue "D47" "0.5679"
Set-TextValue "E47" "  +8.84%  "
Set-TextValue "D48" "125.54"
Set-TextValue "E48" "  +6.12%  "
Set-TextValue "D49" "1.940"
Set-TextValue "E49" "  +7.51%  "
Set-TextValue "D50" "0.06805"
Set-TextValue "E50" "  +5.22%  "
Set-TextValue "D51" "72.55"
Set-TextValue "E51" "  +7.87%  "
